$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2
Set-TextCell "D2" "28.030.92"

# Row 3
Set-TextCell "D3" "1.830.21"
$ws.Range("E3").Value = "  -1.16%  "

# Row 4
$ws.Range("E4").Value = "  -0.26%  "

# Row 5
Set-TextCell "D5" "326.94"
$ws.Range("E5").Value = "  -2.55%  "

# Row 6
Set-TextCell "D6" "1.001"
$ws.Range("E6").Value = "  -0.23%  "

# Row 7
Set-TextCell "D7" "0.4614"
$ws.Range("E7").Value = "  -0.61%  "

# Row 8
Set-TextCell "D8" "0.3865"
$ws.Range("E8").Value = "  -1.62%  "

# Row 9
Set-TextCell "D9" "0.07869"
$ws.Range("E9").Value = "  -0.60%  "

# Row 10
Set-TextCell "D10" "0.9584"
$ws.Range("E10").Value = "  -2.53%  "

# Row 11
Set-TextCell "D11" "21.88"

# Row 12
Set-TextCell "D12" "1.855.26"
$ws.Range("E12").Value = "  +0.25%  "

# Row 13
Set-TextCell "D13" "5.656"
$ws.Range("E13").Value = "  -3.09%  "

# Row 14
Set-TextCell "D14" "6.886"
$ws.Range("E14").Value = "  -1.80%  "

# Row 15
Set-TextCell "D15" "0.06767"
$ws.Range("E15").Value = "  +0.00%  "

# Row 16
$ws.Range("E16").Value = "  -0.33%  "

# Row 17
Set-TextCell "D17" "86.71"
$ws.Range("E17").Value = "  -1.02%  "

# Row 18
Set-TextCell "D18" "0.000009934"
$ws.Range("E18").Value = "  -2.01%  "

# Row 19
Set-TextCell "D19" "16.66"
$ws.Range("E19").Value = "  -2.29%  "

# Row 20
$ws.Range("E20").Value = "  -0.12%  "

# Row 21
Set-TextCell "D21" "28.047.21"
$ws.Range("E21").Value = "  -2.34%  "

# Row 22
Set-TextCell "D22" "5.304"
$ws.Range("E22").Value = "  -2.07%  "

# Row 23
$ws.Range("E23").Value = "  -2.93%  "

# Row 24
Set-TextCell "D24" "2.102"
$ws.Range("E24").Value = "  -1.42%  "

# Row 25
Set-TextCell "D25" "2.095.68"
$ws.Range("E25").Value = "  +1.32%  "

# Row 26
Set-TextCell "D26" "153.73"
$ws.Range("E26").Value = "  +0.17%  "

# Row 27
Set-TextCell "D27" "19.21"
$ws.Range("E27").Value = "  -1.03%  "

# Row 28
Set-TextCell "D28" "5.731"
$ws.Range("E28").Value = "  -9.28%  "

# Row 29
Set-TextCell "D29" "1.972"
$ws.Range("E29").Value = "  -2.72%  "

# Row 30
Set-TextCell "D30" "117.29"
$ws.Range("E30").Value = "  +0.13%  "

# Row 31
Set-TextCell "D31" "0.9361"
$ws.Range("E31").Value = "  -4.53%  "

# Row 32
$ws.Range("E32").Value = "  -1.90%  "

# Row 33
Set-TextCell "D33" "5.294"
$ws.Range("E33").Value = "  -1.91%  "

# Row 34
Set-TextCell "D34" "1.317"
$ws.Range("E34").Value = "  -2.59%  "

# Row 35
Set-TextCell "D35" "3.318"
$ws.Range("E35").Value = "  -4.93%  "

# Row 36
Set-TextCell "D36" "0.05871"
$ws.Range("E36").Value = "  -4.25%  "

# Row 37
$ws.Range("E37").Value = "  -2.56%  "

# Row 38
Set-TextCell "D38" "1.140"
$ws.Range("E38").Value = "  -2.16%  "

# Row 39
Set-TextCell "D39" "7.744"
$ws.Range("E39").Value = "  +1.43%  "

# Row 40
$ws.Range("E40").Value = "  -2.56%  "

# Row 41
Set-TextCell "D41" "9.912"
$ws.Range("E41").Value = "  -2.16%  "

# Row 42
Set-TextCell "D42" "0.1761"
$ws.Range("E42").Value = "  -1.55%  "

# Row 43
Set-TextCell "D43" "1.231"
$ws.Range("E43").Value = "  +0.40%  "

# Row 44
Set-TextCell "D44" "11.61"
$ws.Range("E44").Value = "  -2.36%  "

# Row 45
Set-TextCell "D45" "0.5261"
$ws.Range("E45").Value = "  -2.74%  "

# Row 46
Set-TextCell "D46" "0.07011"
$ws.Range("E46").Value = "  -1.83%  "

# Row 47
Set-TextCell "D47" "2.143"
$ws.Range("E47").Value = "  -10.36%  "

# Row 48
$ws.Range("E48").Value = "  -4.75%  "

# Row 49
Set-TextCell "D49" "112.85"
$ws.Range("E49").Value = "  -2.64%  "

# Row 50
$ws.Range("E50").Value = "  -0.31%  "

# Row 51
$ws.Range("E51").Value = "  -0.41%  "
